# Apply the master-machine_master_h.xlsx edit:
#  - Reformat MAC addresses (colon-separated lowercase -> dash-separated uppercase)
#  - Add 9 new machine rows (Machine 21 - Machine 29)
#  - Widen column C to fit the new mac_address values
#  - Update the sheet selection to match the post-edit state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full data table (rows 2-30, columns A(id) .. K(eff_dtimes)) as a 2D array
$data = New-Object 'object[,]' 29,11
$data[0,0] = 10001
$data[0,1] = "Machine 1"
$data[0,2] = "8C-16-45-5A-5D-0D"
$data[0,3] = "NM5328114630"
$data[0,4] = "192.168.0.150"
$data[0,5] = 1001
$data[0,6] = "eng"
$data[0,7] = $true
$data[0,8] = "superadmin"
$data[0,9] = "now()"
$data[0,10] = "now()"
$data[1,0] = 10002
$data[1,1] = "Machine 2"
$data[1,2] = "8C-16-45-88-E1-0D"
$data[1,3] = "WY2132605316"
$data[1,4] = "192.168.0.133"
$data[1,5] = 1001
$data[1,6] = "eng"
$data[1,7] = $true
$data[1,8] = "superadmin"
$data[1,9] = "now()"
$data[1,10] = "now()"
$data[2,0] = 10003
$data[2,1] = "Machine 3"
$data[2,2] = "00-FF-D3-E3-9A-27"
$data[2,3] = "CM6384145127"
$data[2,4] = "192.168.0.161"
$data[2,5] = 1001
$data[2,6] = "eng"
$data[2,7] = $true
$data[2,8] = "superadmin"
$data[2,9] = "now()"
$data[2,10] = "now()"
$data[3,0] = 10004
$data[3,1] = "Machine 4"
$data[3,2] = "8C-16-45-5A-62-41"
$data[3,3] = "NT894252578"
$data[3,4] = "192.168.0.259"
$data[3,5] = 1001
$data[3,6] = "eng"
$data[3,7] = $true
$data[3,8] = "superadmin"
$data[3,9] = "now()"
$data[3,10] = "now()"
$data[4,0] = 10005
$data[4,1] = "Machine 5"
$data[4,2] = "E8-6A-64-1D-75-E4"
$data[4,3] = "YM866672706"
$data[4,4] = "192.168.0.119"
$data[4,5] = 1001
$data[4,6] = "eng"
$data[4,7] = $true
$data[4,8] = "superadmin"
$data[4,9] = "now()"
$data[4,10] = "now()"
$data[5,0] = 10006
$data[5,1] = "Machine 6"
$data[5,2] = "8C-16-45-FA-94-B7"
$data[5,3] = "WT6501645780"
$data[5,4] = "192.168.0.177"
$data[5,5] = 1001
$data[5,6] = "eng"
$data[5,7] = $true
$data[5,8] = "superadmin"
$data[5,9] = "now()"
$data[5,10] = "now()"
$data[6,0] = 10007
$data[6,1] = "Machine 7"
$data[6,2] = "8C-16-45-1A-0F-62"
$data[6,3] = "LK8186452621"
$data[6,4] = "192.168.0.227"
$data[6,5] = 1001
$data[6,6] = "eng"
$data[6,7] = $true
$data[6,8] = "superadmin"
$data[6,9] = "now()"
$data[6,10] = "now()"
$data[7,0] = 10008
$data[7,1] = "Machine 8"
$data[7,2] = "E8-6A-64-1C-52-6E"
$data[7,3] = "NR3264783870"
$data[7,4] = "192.168.0.207"
$data[7,5] = 1001
$data[7,6] = "eng"
$data[7,7] = $true
$data[7,8] = "superadmin"
$data[7,9] = "now()"
$data[7,10] = "now()"
$data[8,0] = 10009
$data[8,1] = "Machine 9"
$data[8,2] = "48-51-B7-10-35-A6"
$data[8,3] = "RW437027336"
$data[8,4] = "192.168.0.220"
$data[8,5] = 1001
$data[8,6] = "eng"
$data[8,7] = $true
$data[8,8] = "superadmin"
$data[8,9] = "now()"
$data[8,10] = "now()"
$data[9,0] = 10010
$data[9,1] = "Machine 10"
$data[9,2] = "8C-16-45-38-F3-F3"
$data[9,3] = "SI158158531"
$data[9,4] = "192.168.0.242"
$data[9,5] = 1001
$data[9,6] = "eng"
$data[9,7] = $true
$data[9,8] = "superadmin"
$data[9,9] = "now()"
$data[9,10] = "now()"
$data[10,0] = 10011
$data[10,1] = "Machine 11"
$data[10,2] = "D4-3D-7E-58-CC-45"
$data[10,3] = "XF3416823469"
$data[10,4] = "192.168.0.173"
$data[10,5] = 1001
$data[10,6] = "eng"
$data[10,7] = $true
$data[10,8] = "superadmin"
$data[10,9] = "now()"
$data[10,10] = "now()"
$data[11,0] = 10012
$data[11,1] = "Machine 12"
$data[11,2] = "8C-16-45-5A-5D-96"
$data[11,3] = "BW4524978011"
$data[11,4] = "192.168.0.203"
$data[11,5] = 1001
$data[11,6] = "eng"
$data[11,7] = $true
$data[11,8] = "superadmin"
$data[11,9] = "now()"
$data[11,10] = "now()"
$data[12,0] = 10013
$data[12,1] = "Machine 13"
$data[12,2] = "8C-16-45-5A-5D-8E"
$data[12,3] = "DB289579153"
$data[12,4] = "192.168.0.112"
$data[12,5] = 1001
$data[12,6] = "eng"
$data[12,7] = $true
$data[12,8] = "superadmin"
$data[12,9] = "now()"
$data[12,10] = "now()"
$data[13,0] = 10014
$data[13,1] = "Machine 14"
$data[13,2] = "8C-16-45-33-A5-5F"
$data[13,3] = "SI4597903231"
$data[13,4] = "192.168.0.178"
$data[13,5] = 1001
$data[13,6] = "eng"
$data[13,7] = $true
$data[13,8] = "superadmin"
$data[13,9] = "now()"
$data[13,10] = "now()"
$data[14,0] = 10015
$data[14,1] = "Machine 15"
$data[14,2] = "3C-95-09-F9-EA-DF"
$data[14,3] = "TJ7809002958"
$data[14,4] = "192.168.0.267"
$data[14,5] = 1001
$data[14,6] = "eng"
$data[14,7] = $true
$data[14,8] = "superadmin"
$data[14,9] = "now()"
$data[14,10] = "now()"
$data[15,0] = 10016
$data[15,1] = "Machine 16"
$data[15,2] = "8C-16-45-88-E7-0B"
$data[15,3] = "JR6082789079"
$data[15,4] = "192.168.0.149"
$data[15,5] = 1001
$data[15,6] = "eng"
$data[15,7] = $true
$data[15,8] = "superadmin"
$data[15,9] = "now()"
$data[15,10] = "now()"
$data[16,0] = 10017
$data[16,1] = "Machine 17"
$data[16,2] = "B4-69-21-5A-DB-C4"
$data[16,3] = "SA3722889241"
$data[16,4] = "192.168.0.127"
$data[16,5] = 1001
$data[16,6] = "eng"
$data[16,7] = $true
$data[16,8] = "superadmin"
$data[16,9] = "now()"
$data[16,10] = "now()"
$data[17,0] = 10018
$data[17,1] = "Machine 18"
$data[17,2] = "E8-6A-64-1D-48-B7"
$data[17,3] = "RR2683722548"
$data[17,4] = "192.168.0.248"
$data[17,5] = 1001
$data[17,6] = "eng"
$data[17,7] = $true
$data[17,8] = "superadmin"
$data[17,9] = "now()"
$data[17,10] = "now()"
$data[18,0] = 10019
$data[18,1] = "Machine 19"
$data[18,2] = "8C-16-45-59-69-09 "
$data[18,3] = "PO6528391346"
$data[18,4] = "192.168.0.121"
$data[18,5] = 1001
$data[18,6] = "eng"
$data[18,7] = $true
$data[18,8] = "superadmin"
$data[18,9] = "now()"
$data[18,10] = "now()"
$data[19,0] = 10020
$data[19,1] = "Machine 20"
$data[19,2] = "98-E7-F4-30-16-5A "
$data[19,3] = "FB5962911652"
$data[19,4] = "192.168.0.215"
$data[19,5] = 1001
$data[19,6] = "eng"
$data[19,7] = $true
$data[19,8] = "superadmin"
$data[19,9] = "now()"
$data[19,10] = "now()"
$data[20,0] = 10021
$data[20,1] = "Machine 21"
$data[20,2] = "38-BA-F8-53-C7-8F"
$data[20,3] = "FB5962911653"
$data[20,4] = "192.168.0.874"
$data[20,5] = 1001
$data[20,6] = "eng"
$data[20,7] = $true
$data[20,8] = "superadmin"
$data[20,9] = "now()"
$data[20,10] = "now()"
$data[21,0] = 10022
$data[21,1] = "Machine 22"
$data[21,2] = "E8-6A-64-1C-58-C2"
$data[21,3] = "FB5962911654"
$data[21,4] = "192.168.0.721"
$data[21,5] = 1001
$data[21,6] = "eng"
$data[21,7] = $true
$data[21,8] = "superadmin"
$data[21,9] = "now()"
$data[21,10] = "now()"
$data[22,0] = 10023
$data[22,1] = "Machine 23"
$data[22,2] = "E4-A4-71-CE-BA-93"
$data[22,3] = "FB5962911655"
$data[22,4] = "192.168.0.841"
$data[22,5] = 1001
$data[22,6] = "eng"
$data[22,7] = $true
$data[22,8] = "superadmin"
$data[22,9] = "now()"
$data[22,10] = "now()"
$data[23,0] = 10024
$data[23,1] = "Machine 24"
$data[23,2] = "54-E1-AD-EA-30-C9"
$data[23,3] = "FB5962911656"
$data[23,4] = "192.168.0.186"
$data[23,5] = 1001
$data[23,6] = "eng"
$data[23,7] = $true
$data[23,8] = "superadmin"
$data[23,9] = "now()"
$data[23,10] = "now()"
$data[24,0] = 10025
$data[24,1] = "Machine 25"
$data[24,2] = "8C-16-45-65-DD-40"
$data[24,3] = "FB5962911657"
$data[24,4] = "192.168.0.627"
$data[24,5] = 1001
$data[24,6] = "eng"
$data[24,7] = $true
$data[24,8] = "superadmin"
$data[24,9] = "now()"
$data[24,10] = "now()"
$data[25,0] = 10026
$data[25,1] = "Machine 26"
$data[25,2] = "58-20-B1-D6-C3-BE"
$data[25,3] = "FB5962911658"
$data[25,4] = "192.168.0.879"
$data[25,5] = 1001
$data[25,6] = "eng"
$data[25,7] = $true
$data[25,8] = "superadmin"
$data[25,9] = "now()"
$data[25,10] = "now()"
$data[26,0] = 10027
$data[26,1] = "Machine 27"
$data[26,2] = "8C-16-45-38-F0-25"
$data[26,3] = "FB5962911659"
$data[26,4] = "192.168.0.628"
$data[26,5] = 1001
$data[26,6] = "eng"
$data[26,7] = $true
$data[26,8] = "superadmin"
$data[26,9] = "now()"
$data[26,10] = "now()"
$data[27,0] = 10028
$data[27,1] = "Machine 28"
$data[27,2] = "6C-88-14-AC-EF-55"
$data[27,3] = "FB5962911661"
$data[27,4] = "192.168.0.306"
$data[27,5] = 1001
$data[27,6] = "eng"
$data[27,7] = $true
$data[27,8] = "superadmin"
$data[27,9] = "now()"
$data[27,10] = "now()"
$data[28,0] = 10029
$data[28,1] = "Machine 29"
$data[28,2] = "3C-6A-A7-C0-DF-27"
$data[28,3] = "FB5962911662"
$data[28,4] = "192.168.0.355"
$data[28,5] = 1001
$data[28,6] = "eng"
$data[28,7] = $true
$data[28,8] = "superadmin"
$data[28,9] = "now()"
$data[28,10] = "now()"

# Write the rebuilt table back into the sheet (replaces existing 20 rows and adds 9 new ones)
$ws.Range("A2:K30").Value = $data

# Column C (mac_address) needs to be a bit wider to fit the new values
$ws.Columns.Item(3).ColumnWidth = 16.14

# Select the row below the data, spanning the full row width (mirrors the saved selection state)
$ws.Range("A31:XFD1048576").Select() | Out-Null

